$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 48
$prevRow = 47

# Row 48 repeats the same Date/Weekday/Week as row 47 (only the time and the
# resale numbers changed), so clone row 47's formatting/types down into row 48
# first (this keeps the date-like text cells, e.g. "2024-01-11" and "01", as
# plain text instead of Excel reinterpreting them as dates/numbers) and then
# overwrite the cells that actually differ.
$ws.Range("A$prevRow`:T$prevRow").Copy() | Out-Null
$ws.Range("A$row`:T$row").PasteSpecial() | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 2).Value = "19:02:26"

$ws.Cells.Item($row, 5).Value = 139554
$ws.Cells.Item($row, 6).Value = 142813
$ws.Cells.Item($row, 7).Value = 171893
$ws.Cells.Item($row, 8).Value = 148245
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 119573
$ws.Cells.Item($row, 11).Value = 224928
$ws.Cells.Item($row, 12).Value = 252345
$ws.Cells.Item($row, 13).Value = 185193
$ws.Cells.Item($row, 14).Value = 110449
$ws.Cells.Item($row, 15).Value = 40789
$ws.Cells.Item($row, 16).Value = 30895
$ws.Cells.Item($row, 17).Value = 72918
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42370
$ws.Cells.Item($row, 20).Value = -1
